# Applies the "update categories, insert empty field complete" edit to the
# VSIG Trial Balance worksheet:
#   - fills in the report header (company name / title / period) in E2, E7, E8
#     and clears E3:E6 explicitly
#   - adds a new "Category" column (H) value for each account row, grouping
#     the individual GL accounts into higher level reporting categories
#   - adds the cross-check formula in H56 (=G56-F56)
#   - moves the active selection to D24 (matches author's last edit position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Report header block (rows 2-8) ---------------------------------------
$ws.Range("E2").Value = "VSIG Pte. Ltd."
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = "Trial Balance"
$ws.Range("E8").Value = "December 2015"

# ---- New "Category" column (H) --------------------------------------------
$categories = @{
    11 = "Bank Balances"
    12 = "Bank Balances"
    13 = "Bank Balances"
    14 = "Trade Receivables"
    15 = "Trade Receivables"
    16 = "Plant and Equipment"
    17 = "Plant and Equipment"
    18 = "Plant and Equipment"
    19 = "Plant and Equipment"
    20 = "Deposits"
    21 = "Prepayments"
    22 = "Trade Payables"
    23 = "Trade Payables"
    26 = "GST Payables"
    27 = "Accruals"
    28 = "Amount owing to a Shareholder"
    29 = "Income Tax Payables"
    30 = "Share Capital"
    31 = "Retained Profits"
    32 = "Revenue"
    33 = "Cost of Sales"
    34 = "Accounting Fee"
    35 = "Administrative Expenses"
    36 = "Bank Charges"
    37 = "Compilation Fee"
    38 = "Depreciation"
    39 = "Entertainment"
    40 = "Freight Charges"
    41 = "Internet Expenses"
    42 = "Late Penalty"
    43 = "Nominee Director Fee"
    44 = "Office Supplies"
    45 = "Postage and Courier"
    46 = "Professional Fee"
    47 = "Secretarial Fee"
    48 = "Taxation Fee"
    49 = "Telephone Expenses"
    50 = "Salaries"
    51 = "Skill Development Levy & SINDA"
    52 = "Exchange Gain - Trade"
    53 = "Exchange Gain - Non-trade"
    54 = "Income Tax Expense"
}

foreach ($row in $categories.Keys) {
    $ws.Cells.Item($row, 8).Value = $categories[$row]
}

# ---- Cross-check formula ----------------------------------------------------
$ws.Range("H56").Formula = "=G56-F56"

# ---- Restore author's last selection ---------------------------------------
$ws.Range("D24").Select() | Out-Null
